$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data rows 2 and 3 (keep A and I and J the same, B/D/E/G/H swap,
# F normalized from -0 to 0) to match the new schedule data.
$ws.Range("B2").Value = 7
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 11

$ws.Range("B3").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16
